# Generate Report for Handback
# Adds a new handback record (file 30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md, now
# "in sync with en-US") as row 3 on the Overview, zh-cn and de-de sheets, and
# grows each sheet's table to include the new row.

$wb = $excel.ActiveWorkbook

$fileName   = "30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md"
$pathName   = "e2e\30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"

$overviewDate = "2016-10-27 07:58:14"

$zhTargetFile = "30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.05175b977d06e199acb47d3e01dc0f8283e0fcdf.zh-cn.xlf"
$zhHoDate     = "2016-10-27 07:58:01"
$zhHbDate     = "2016-10-27 07:58:43"

$deTargetFile = "30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.05175b977d06e199acb47d3e01dc0f8283e0fcdf.de-de.xlf"
$deHoDate     = "2016-10-27 07:58:14"
$deHbDate     = "2016-10-27 07:59:01"

$srcPath  = "e2e"
$priority = "ht"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathName
$wsOverview.Range("C3").Value = $ext
$wsOverview.Range("E3").Value = $statusSync
$wsOverview.Range("F3").Value = $statusSync
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f20fbe22e3821466bad32cefe11fd9332e4db404/e2e/30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md", "", "", $pathName)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $fileName
$wsZh.Range("B3").Value = $ext
$wsZh.Range("C3").Value = $statusSync
$wsZh.Range("D3").Value = $srcPath
$wsZh.Range("E3").Value = $priority
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = $zhTargetFile
$wsZh.Range("H3").Value = $zhHoDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = $fileName
$wsZh.Range("J3").Value = $zhTargetFile
$wsZh.Range("K3").Value = $zhHbDate
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f20fbe22e3821466bad32cefe11fd9332e4db404/e2e/30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md", "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c4346497500f4654f614cc4178d21ce0c555c8d0/e2e/30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md", "", "", $fileName)

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $fileName
$wsDe.Range("B3").Value = $ext
$wsDe.Range("C3").Value = $statusSync
$wsDe.Range("D3").Value = $srcPath
$wsDe.Range("E3").Value = $priority
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = $deTargetFile
$wsDe.Range("H3").Value = $deHoDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = $fileName
$wsDe.Range("J3").Value = $deTargetFile
$wsDe.Range("K3").Value = $deHbDate
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f20fbe22e3821466bad32cefe11fd9332e4db404/e2e/30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md", "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bf48de1bf58d909e61a3f82adf147138013e0781/e2e/30d22998-9b4a-4e5a-a0cd-46ba3a8bb5ea.md", "", "", $fileName)

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
